$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(0.01384957621052772, 2.21878651858606, 21.84152086182795, 4.673491292580734, 4.783449836594167, 22)
    3  = @(-0.1108950891496856, 2.52800863553415, 20.81929224357433, 4.562816262307122, 4.674114274267978, 21)
    4  = @(-0.5909552399615555, 1.803240045667984, 9.880028153484647, 3.143251207505478, 3.167399740090089, 20)
    5  = @(-0.1812638745297453, 1.430346199258623, 9.695185292776637, 3.113709249878132, 3.193606854246072, 19)
    6  = @(-0.1687010319474557, 1.931467188986993, 13.19207325907522, 3.632089379279539, 3.733355516202483, 18)
    7  = @(-0.2974795460340417, 1.914100224451346, 11.03882338163619, 3.322472480192453, 3.410971252101159, 17)
    8  = @(-0.1960622915985716, 1.836194656235551, 11.49912085291426, 3.391035365919127, 3.496387534695756, 16)
    9  = @(-0.08754086917198256, 1.885964882007152, 11.84300447804472, 3.441366658472288, 3.561000221788015, 15)
    10 = @(-0.1199187714296267, 2.158918608588491, 14.23208708490037, 3.772543847975842, 3.912975378092231, 14)
    11 = @(-0.08524854279770107, 2.086633922413215, 13.79239774298984, 3.713811753843999, 3.864439325182775, 13)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    $ws.Cells.Item($row, 2).Value = $values[0]
    $ws.Cells.Item($row, 3).Value = $values[1]
    $ws.Cells.Item($row, 4).Value = $values[2]
    $ws.Cells.Item($row, 5).Value = $values[3]
    $ws.Cells.Item($row, 6).Value = $values[4]
    $ws.Cells.Item($row, 7).Value = $values[5]
}
